$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel keeps them as text (matching original inlineStr formatting)
# instead of auto-converting them to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.517.22"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "2.660.75"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "605.05"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "157.90"
$ws.Range("E6").Value = "  +4.64%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("E9").Value = "  +7.92%  "
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "29.65"
$ws.Range("E13").Value = "  +6.05%  "
$ws.Range("E14").Value = "  +15.98%  "
$ws.Range("D15").Value = "3.139.00"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "65.336.74"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "2.650.84"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "12.82"
$ws.Range("E18").Value = "  +4.68%  "
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "360.46"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "7.37"
$ws.Range("E21").Value = "  +5.09%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").Value = "9.56"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  +16.99%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "8.26"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("E30").Value = "  +8.01%  "
$ws.Range("D31").Value = "542.15"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  +4.94%  "
$ws.Range("D35").Value = "6.44"
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("D36").Value = "0.434"
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("D37").Value = "20.66"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("D38").Value = "163.10"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "42.54"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("D43").Value = "166.69"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  +7.79%  "
$ws.Range("E46").Value = "  +5.55%  "
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +4.96%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "19.83"
$ws.Range("E51").Value = "  +2.98%  "

# Revert the temporary text-format styling so the cells keep their original
# (unstyled) appearance, now holding text values.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D51").Style = "Normal"
